# Rename ward codes to include city-name prefixes (wardoslo/wardstavanger/wardbergen/wardtrondheim)
# and add a duplicate set of Oslo ward rows using the new "wardoslo" naming.
#
# Layout before edit (rows 2-49):
#   2-16   Oslo      (year_start=2005)                          ward_code = ward030xxx
#   17-25  Stavanger (year_start=2005)                          ward_code = ward1103xx
#   26-33  Bergen    (year_start=2005, year_end=2019)            ward_code = ward1201xx, ward_code_end = ward4601xx
#   34-41  Bergen    (year_start=2020)                          ward_code = ward4601xx
#   42-45  Trondheim (year_start=2005, year_end=2017)            ward_code = ward1601xx, ward_code_end = ward5001xx
#   46-49  Trondheim (year_start=2018)                           ward_code = ward5001xx
#
# After edit:
#   2-16   Oslo        unchanged
#   17-31  Oslo        NEW duplicate rows, ward_code = wardoslo030xxx
#   32-40  Stavanger   ward_code -> wardstavanger1103xx
#   41-48  Bergen      ward_code -> wardbergen1201xx, ward_code_end -> wardbergen4601xx
#   49-56  Bergen      ward_code -> wardbergen4601xx
#   57-60  Trondheim   ward_code -> wardtrondheim1601xx, ward_code_end -> wardtrondheim5001xx
#   61-64  Trondheim   ward_code -> wardtrondheim5001xx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 15 blank rows right after row 16 (new rows 17..31), pushing everything else down.
$ws.Rows("17:31").Insert()

# 2) Fill the new rows 17-31 as copies of (the now-shifted) rows 2-16, but rename
#    the ward_code (column E) to use the "wardoslo" prefix.
for ($i = 0; $i -lt 15; $i++) {
    $srcRow = 2 + $i
    $dstRow = 17 + $i

    $yearStart = $ws.Cells.Item($srcRow, 1).Value2
    $wardCode  = $ws.Cells.Item($srcRow, 5).Value2
    $wardName  = $ws.Cells.Item($srcRow, 6).Value2
    $municipCode = $ws.Cells.Item($srcRow, 7).Value2
    $municipName = $ws.Cells.Item($srcRow, 8).Value2

    $newWardCode = $wardCode -replace '^ward', 'wardoslo'

    $ws.Cells.Item($dstRow, 1).Value = $yearStart
    $ws.Cells.Item($dstRow, 5).Value = $newWardCode
    $ws.Cells.Item($dstRow, 6).Value = $wardName
    $ws.Cells.Item($dstRow, 7).Value = $municipCode
    $ws.Cells.Item($dstRow, 8).Value = $municipName
}

# 3) Rename the ward codes (and ward_code_end where present) for the existing
#    Stavanger / Bergen / Trondheim rows, now shifted down by 15 rows, i.e. rows 32-64.
$dim = $ws.UsedRange
$lastRow = $dim.Row + $dim.Rows.Count - 1

for ($r = 32; $r -le $lastRow; $r++) {
    $wardCodeCell = $ws.Cells.Item($r, 5)
    $wardCode = $wardCodeCell.Value2
    if ($wardCode -ne $null -and $wardCode -ne "") {
        if ($wardCode -match '^ward110') {
            $wardCodeCell.Value = $wardCode -replace '^ward', 'wardstavanger'
        } elseif ($wardCode -match '^ward120') {
            $wardCodeCell.Value = $wardCode -replace '^ward', 'wardbergen'
        } elseif ($wardCode -match '^ward460') {
            $wardCodeCell.Value = $wardCode -replace '^ward', 'wardbergen'
        } elseif ($wardCode -match '^ward160') {
            $wardCodeCell.Value = $wardCode -replace '^ward', 'wardtrondheim'
        } elseif ($wardCode -match '^ward500') {
            $wardCodeCell.Value = $wardCode -replace '^ward', 'wardtrondheim'
        }
    }

    $wardCodeEndCell = $ws.Cells.Item($r, 3)
    $wardCodeEnd = $wardCodeEndCell.Value2
    if ($wardCodeEnd -ne $null -and $wardCodeEnd -ne "") {
        if ($wardCodeEnd -match '^ward460') {
            $wardCodeEndCell.Value = $wardCodeEnd -replace '^ward', 'wardbergen'
        } elseif ($wardCodeEnd -match '^ward500') {
            $wardCodeEndCell.Value = $wardCodeEnd -replace '^ward', 'wardtrondheim'
        }
    }
}

# 4) Column E width: widen to fit the longer "wardXXXXXXXXXX" codes.
$ws.Columns("E").ColumnWidth = 21

# 5) Selection / view bookkeeping to mirror the saved workbook state.
$ws.Range("G61").Select()
